# Adds 4 new annotation rows (155-158) to the "Tabela1" Excel table on
# sheet "Planilha1", for lesson 64 ("Implementando SmtpEmailService com
# servidor do Google"), mirroring the existing B:G layout
# (Seção | Nome da Seção | Aula | nome aula | abordagem da aula | aprendido).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")
$lo = $ws.ListObjects.Item("Tabela1")

# Seed the 4 new rows with the formatting (cell styles incl. wrap-text/
# alignment) used by the existing last table row, so newly appended cells
# look consistent with the rest of the table.
$ws.Range("B154:G154").Copy()
$ws.Range("B155:G158").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# "abordagem da aula" (column F) for each new row - filled first, matching
# the order the notes were originally authored in.
$ws.Range("F155").Value = "4:22`n5. Serviço de email`n64. Implementando SmtpEmailService com servidor do Google`ninstanciação de MailSender que é uma classe do framework - ele automaticamente pega os dados de e-mail do arquivo application.properties para enviar e-mails"
$ws.Range("F156").Value = "5:18`n5. Serviço de email`n64. Implementando SmtpEmailService com servidor do Google`ncriação da classe SmtpEmailService - responsável por pegar os dados de configs de login e senha do arquivo application.properties e enviar e-mails"
$ws.Range("F157").Value = "7:06`n5. Serviço de email`n64. Implementando SmtpEmailService com servidor do Google`ncriação de @Bean na classe DevConfig para instanciação de SmtpEmailService (envio de e-mail)"
$ws.Range("F158").Value = "8:36`n5. Serviço de email`n64. Implementando SmtpEmailService com servidor do Google`nIMPORTANTE: o envio de e-mail (PELO GOGGLE SMTP) só funcionou para mim utilizando a solução descrita na aula (acessando a aba Segurança>Acesso a app menos seguro>ATIVAR ou pelo link https://myaccount.google.com/lesssecureapps) porém foi necessário também desativar o antivirus temporariamente"

# "nome aula" (column E) - same lesson title repeated on every new row.
$ws.Range("E155").Value = "Implementando SmtpEmailService com servidor do Google"
$ws.Range("E156").Value = "Implementando SmtpEmailService com servidor do Google"
$ws.Range("E157").Value = "Implementando SmtpEmailService com servidor do Google"
$ws.Range("E158").Value = "Implementando SmtpEmailService com servidor do Google"

# "Seção" (B), "Nome da Seção" (C) and "Aula" (D) for each new row.
$ws.Range("B155").Value = 5
$ws.Range("C155").Value = "Serviço de email"
$ws.Range("D155").Value = 64
$ws.Range("B156").Value = 5
$ws.Range("C156").Value = "Serviço de email"
$ws.Range("D156").Value = 64
$ws.Range("B157").Value = 5
$ws.Range("C157").Value = "Serviço de email"
$ws.Range("D157").Value = 64
$ws.Range("B158").Value = 5
$ws.Range("C158").Value = "Serviço de email"
$ws.Range("D158").Value = 64

# "aprendido" (column G) - only the first new row carries the little
# blank-line marker used elsewhere in the sheet; the other 3 stay empty.
$ws.Range("G155").Value = "`n`n`n`n`n`n"

# Grow the table (ListObject) so the new rows participate in the filter/
# banding, matching ref="B1:G158" / autoFilter ref="B1:G158".
$lo.Resize($ws.Range("B1:G158"))

# Explicit row heights (wrap-text autosize equivalents) for the new rows.
$ws.Rows.Item(155).RowHeight = 105
$ws.Rows.Item(156).RowHeight = 90
$ws.Rows.Item(157).RowHeight = 75
$ws.Rows.Item(158).RowHeight = 120

# Scroll position / active selection, matching the saved view state.
$ws.Application.ActiveWindow.ScrollRow = 144
$ws.Range("E155").Select()

# Window height tweak recorded in the workbook view.
$excel.ActiveWindow.Height = 11010
